$wb = $excel.ActiveWorkbook

# Rename the original sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Mentioned_in_text"

# Add the new sheet right after it, and populate it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Extra_on_github"

$ws2.Range("A1").Value = "Some extra files that do not need to be mentioned in the text but that should anyway be provided in the github:"
$ws2.Range("A3").Value = "scripts..."
$ws2.Range("A4").Value = "list of marker genes used for each species in the metacell pipeline.... The gene codes would be present in the metacell script for each species, but perhaps provide also a fasta file?.... Also the description of which categorise... See google sheets files..."
$ws2.Range("A5").Value = "for the metacell pipelines: in theory I could provide also all the files necessary for people to run the scripts..... Separate directory per species..."

# Update selection on the first sheet.
$ws1.Range("C13").Select() | Out-Null

# Update selection on the new sheet (it becomes the active/selected sheet).
$ws2.Range("E12").Select() | Out-Null
